# Update ChanjoKe FHIR IG metadata and element values

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "https://intellisoft-consulting.github.io/igs/ChanjoKe-FHIR-IG/StructureDefinition/kenya-counties-extension"
$meta.Range("B8").Value = "2024-08-27T20:30:12+00:00"
$meta.Range("B9").Value = "Intellisoft Consulting Ltd"
$meta.Range("B10").Value = "Intellisoft Consulting Ltd (https://www.intellisoftkenya.com/, info[at]intellisoftkenya.com)"
$meta.Range("B21").Value = "element:Patient"

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")

# R5 (Fixed Value) shared the same text as Metadata!B2 (URL) in the original
# workbook, so it must be kept in sync with the new StructureDefinition URL.
$elem.Range("R5").Value = "https://intellisoft-consulting.github.io/igs/ChanjoKe-FHIR-IG/StructureDefinition/kenya-counties-extension"

$elem.Range("Z6").Value = "https://intellisoft-consulting.github.io/igs/ChanjoKe-FHIR-IG/ValueSet/kenyaCountiesVS"

# Column Z (Binding Value Set) widens to fit the new, longer URL -
# target stored width ~80.43 characters.
$elem.Columns.Item(26).ColumnWidth = 79.6
